# Apply edits described in the commit:
# "removed C8 and C10 from all variants as tests showed that without
#  these the SPI can be clocked a little faster"
#
# Concretely, for the BOM row that used to list designators
# "C2, C3, C4, C5, C8, C10" (row 3 of the sheet), C8 and C10 are removed
# from the designator list (becoming "C2, C3, C4, C5") and the Quantity
# of that row drops from 6 to 4 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Quantity (A3) 6 -> 4
$ws.Range("A3").Value2 = 4

# Row 3: Designator (B3) "C2, C3, C4, C5, C8, C10" -> "C2, C3, C4, C5"
$ws.Range("B3").Value2 = "C2, C3, C4, C5"

# Minor formatting touch-up that came along with the edit/save
$ws.Columns("D").ColumnWidth = 19.5

# Leave the cursor/selection on the edited row
$ws.Range("A3").Select()

Write-Host "Edit applied"
